$d = $word.ActiveDocument

function Set-BoldAndText($OldText, $NewText, $ReplaceText) {
    # Step 1: locate the run and make it bold (w:b w:val="0" -> w:b/>)
    $findRange = $d.Content
    $find = $findRange.Find
    $find.ClearFormatting()
    $found = $find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $findRange.Font.Bold = 1
    }

    # Step 2: replace the text itself, if requested
    if ($ReplaceText) {
        $replaceRange = $d.Content
        $rfind = $replaceRange.Find
        $rfind.ClearFormatting()
        $rfind.Replacement.ClearFormatting()
        [void]$rfind.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, $NewText, 2)
    }
}

Set-BoldAndText "매출(`$1,000)" "수익(`$K)" $true
Set-BoldAndText "판매 제품 원가(`$1,000)" "판매 제품 원가(`$K)" $true
Set-BoldAndText "총 이익률(%)" "" $false
Set-BoldAndText "영업비용(`$1,000)" "영업 비용(`$K)" $true
Set-BoldAndText "EBITDA(`$1,000)" "EBITDA(`$K)" $true
Set-BoldAndText "이자비용(`$1,000)" "이자 비용(`$K)" $true
Set-BoldAndText "세전 이익(`$1,000)" "세전 이익(`$K)" $true
Set-BoldAndText "순수입(`$1,000)" "순이익(`$K)" $true
Set-BoldAndText "자산 총액(`$1,000)" "총 자산(`$K)" $true
Set-BoldAndText "부채 총액(`$1,000)" "총 부채(`$K)" $true
Set-BoldAndText "자기 자본(`$1,000)" "주주 지분(`$K)" $true

Write-Host "Done."
